$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that were removed in the diff
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("AN7").ClearContents()
$ws.Range("AP7").ClearContents()

# Set / update cell values per diff
$ws.Range("K2").Value = -25.278
$ws.Range("U2").Value = 8.488
$ws.Range("V2").Value = 0.02310288513881328
$ws.Range("W2").Value = -0.09807692307692308
$ws.Range("X2").Value = 0.07737508688094547
$ws.Range("Y2").Value = -0.1754520099578685
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = -0.08111640781177917
$ws.Range("AB2").Value = 0.07664805266163008
$ws.Range("AC2").Value = -0.1569858619190282
$ws.Range("AD2").Value = 6.169
$ws.Range("AF2").Value = 6.169
$ws.Range("AG2").Value = -2.318999999999999
$ws.Range("AH2").Value = 0.01651368288053881
$ws.Range("AI2").Value = 0.05834176604658641
$ws.Range("AJ2").Value = -0.006352015032280505
$ws.Range("AK2").Value = -0.02384551315667704
$ws.Range("AL2").Value = 0.9350000000000001
$ws.Range("AM2").Value = 0.914
$ws.Range("AN2").Value = -0.9382509505703422
$ws.Range("AO2").Value = -8.126203208556149
$ws.Range("AP2").Value = 0.352699619771863
$ws.Range("AQ2").Value = -8.312910284463895
$ws.Range("B3").Value = "Talon Metals Corp. (TSX:TLO)"
$ws.Range("K3").Value = -0.838
$ws.Range("O3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 5.36
$ws.Range("V3").Value = 0.02245496439044826
$ws.Range("W3").Value = -0.02004784688995215
$ws.Range("X3").Value = 0.07586945410724899
$ws.Range("Y3").Value = -0.09591730099720114
$ws.Range("AA3").Value = -0.01605600238308014
$ws.Range("AB3").Value = 0.07586945410724899
$ws.Range("AC3").Value = -0.09192545649032913
$ws.Range("AG3").Value = -5.36
$ws.Range("AJ3").Value = -0.02297077226364961
$ws.Range("AK3").Value = -0.123959296947271
$ws.Range("AM3").Value = -0.02
$ws.Range("AN3").Value = -0
$ws.Range("AP3").Value = 10.20952380952381
$ws.Range("AQ3").Value = 26.95
$ws.Range("K4").Value = -28.8
$ws.Range("U4").Value = 0.21
$ws.Range("V4").Value = 0.004046242774566474
$ws.Range("W4").Value = -0.7559055118110236
$ws.Range("X4").Value = 0.08010757595663927
$ws.Range("Y4").Value = -0.8360130877676629
$ws.Range("AA4").Value = -0.06519507186858316
$ws.Range("AB4").Value = 0.07771328146294879
$ws.Range("AC4").Value = -0.142908353331532
$ws.Range("AD4").Value = 3.98
$ws.Range("AF4").Value = 3.98
$ws.Range("AG4").Value = 3.77
$ws.Range("AH4").Value = 0.07122405153901218
$ws.Range("AI4").Value = 0.3601809954751131
$ws.Range("AJ4").Value = 0.06772049577869589
$ws.Range("AK4").Value = 0.3477859778597786
$ws.Range("AL4").Value = 0.377
$ws.Range("AM4").Value = 0.377
$ws.Range("AN4").Value = -1.604838709677419
$ws.Range("AO4").Value = -6.737400530503979
$ws.Range("AP4").Value = -1.520161290322581
$ws.Range("AQ4").Value = -6.737400530503979
$ws.Range("B5").Value = "Amur Minerals Corporation (AIM:AMC)"
$ws.Range("K5").Value = -2.53
$ws.Range("U5").Value = 0.831
$ws.Range("V5").Value = 0.02760797342192691
$ws.Range("W5").Value = -0.09999999999999999
$ws.Range("X5").Value = 0.07586945410724899
$ws.Range("Y5").Value = -0.175869454107249
$ws.Range("AA5").Value = -0.08111640781177917
$ws.Range("AB5").Value = 0.07586945410724899
$ws.Range("AC5").Value = -0.1569858619190282
$ws.Range("AG5").Value = -0.831
$ws.Range("AJ5").Value = -0.02839181386449827
$ws.Range("AK5").Value = -0.03438288716951467
$ws.Range("AL5").Value = 0.462
$ws.Range("AM5").Value = 0.461
$ws.Range("AO5").Value = -4.567099567099566
$ws.Range("AQ5").Value = -4.57700650759219
$ws.Range("B6").Value = "Premier African Minerals Limited (AIM:PREM)"
$ws.Range("K6").Value = 7.91
$ws.Range("O6").Value = -0
$ws.Range("R6").Value = -0
$ws.Range("U6").Value = 0.037
$ws.Range("V6").Value = 0.002936507936507936
$ws.Range("W6").Value = 0.4598837209302326
$ws.Range("X6").Value = 0.08139605018723281
$ws.Range("Y6").Value = 0.3784876707429998
$ws.Range("AA6").Value = -0.08531710518465485
$ws.Range("AB6").Value = 0.07936440698839346
$ws.Range("AC6").Value = -0.1646815121730483
$ws.Range("AD6").Value = 1.26
$ws.Range("AF6").Value = 1.26
$ws.Range("AG6").Value = 1.223
$ws.Range("AH6").Value = 0.09090909090909091
$ws.Range("AI6").Value = 0.1891891891891892
$ws.Range("AJ6").Value = 0.08847572885770094
$ws.Range("AK6").Value = 0.1846595198550506
$ws.Range("AL6").Value = 0.044
$ws.Range("AM6").Value = 0.044
$ws.Range("AN6").Value = -0.3529411764705883
$ws.Range("AO6").Value = -32.5
$ws.Range("AP6").Value = -0.342577030812325
$ws.Range("AQ6").Value = -32.5
$ws.Range("B7").Value = "Phoenix Copper Limited (LSE:PXC)"
$ws.Range("K7").Value = -1.02
$ws.Range("U7").Value = 2.05
$ws.Range("V7").Value = 0.06011730205278592
$ws.Range("W7").Value = -0.09807692307692308
$ws.Range("X7").Value = 0.07737508688094547
$ws.Range("Y7").Value = -0.1754520099578685
$ws.Range("Z7").Value = 0
$ws.Range("AA7").Value = -0.09708449028163427
$ws.Range("AB7").Value = 0.07664805266163008
$ws.Range("AC7").Value = -0.1737325429432643
$ws.Range("AD7").Value = 0.929
$ws.Range("AF7").Value = 0.929
$ws.Range("AG7").Value = -1.121
$ws.Range("AH7").Value = 0.02652088269719375
$ws.Range("AI7").Value = 0.06438422621110264
$ws.Range("AJ7").Value = -0.03399132781466994
$ws.Range("AK7").Value = -0.09055658776960981
$ws.Range("AL7").Value = 0.052
$ws.Range("AM7").Value = 0.052
$ws.Range("AO7").Value = -18.82692307692308
$ws.Range("AQ7").Value = -18.82692307692308
